$d = $word.ActiveDocument

$d.Content.Find.Execute("631×5=3155", $true, $false, $false, $false, $false, $true, 1, $false, "949×3=2847", 2) | Out-Null
$d.Content.Find.Execute("989×4=3956", $true, $false, $false, $false, $false, $true, 1, $false, "780×3=2340", 2) | Out-Null
$d.Content.Find.Execute("300×9=2700", $true, $false, $false, $false, $false, $true, 1, $false, "294×4=1176", 2) | Out-Null
$d.Content.Find.Execute("115×4=460", $true, $false, $false, $false, $false, $true, 1, $false, "329×8=2632", 2) | Out-Null
$d.Content.Find.Execute("678×6=4068", $true, $false, $false, $false, $false, $true, 1, $false, "857×6=5142", 2) | Out-Null
$d.Content.Find.Execute("269×6=1614", $true, $false, $false, $false, $false, $true, 1, $false, "826×8=6608", 2) | Out-Null
$d.Content.Find.Execute("574×8=4592", $true, $false, $false, $false, $false, $true, 1, $false, "618×9=5562", 2) | Out-Null
$d.Content.Find.Execute("932×7=6524", $true, $false, $false, $false, $false, $true, 1, $false, "601×4=2404", 2) | Out-Null
$d.Content.Find.Execute("944×7=6608", $true, $false, $false, $false, $false, $true, 1, $false, "800×4=3200", 2) | Out-Null
$d.Content.Find.Execute("400×9=3600", $true, $false, $false, $false, $false, $true, 1, $false, "857×6=5142", 2) | Out-Null
$d.Content.Find.Execute("946×7=6622", $true, $false, $false, $false, $false, $true, 1, $false, "972×9=8748", 2) | Out-Null
$d.Content.Find.Execute("531×6=3186", $true, $false, $false, $false, $false, $true, 1, $false, "658×2=1316", 2) | Out-Null
$d.Content.Find.Execute("461×2=922", $true, $false, $false, $false, $false, $true, 1, $false, "742×3=2226", 2) | Out-Null
$d.Content.Find.Execute("536×9=4824", $true, $false, $false, $false, $false, $true, 1, $false, "498×3=1494", 2) | Out-Null
$d.Content.Find.Execute("718×5=3590", $true, $false, $false, $false, $false, $true, 1, $false, "513×3=1539", 2) | Out-Null
$d.Content.Find.Execute("929×2=1858", $true, $false, $false, $false, $false, $true, 1, $false, "554×9=4986", 2) | Out-Null
$d.Content.Find.Execute("388×8=3104", $true, $false, $false, $false, $false, $true, 1, $false, "399×7=2793", 2) | Out-Null
$d.Content.Find.Execute("713×2=1426", $true, $false, $false, $false, $false, $true, 1, $false, "919×7=6433", 2) | Out-Null
$d.Content.Find.Execute("522×4=2088", $true, $false, $false, $false, $false, $true, 1, $false, "565×2=1130", 2) | Out-Null
$d.Content.Find.Execute("890×9=8010", $true, $false, $false, $false, $false, $true, 1, $false, "258×5=1290", 2) | Out-Null
$d.Content.Find.Execute("362×2=724", $true, $false, $false, $false, $false, $true, 1, $false, "273×8=2184", 2) | Out-Null
$d.Content.Find.Execute("366×4=1464", $true, $false, $false, $false, $false, $true, 1, $false, "964×7=6748", 2) | Out-Null
$d.Content.Find.Execute("637×5=3185", $true, $false, $false, $false, $false, $true, 1, $false, "817×8=6536", 2) | Out-Null
$d.Content.Find.Execute("498×2=996", $true, $false, $false, $false, $false, $true, 1, $false, "672×8=5376", 2) | Out-Null
$d.Content.Find.Execute("800×2=1600", $true, $false, $false, $false, $false, $true, 1, $false, "457×3=1371", 2) | Out-Null
